# Auto-generated Excel COM-interop script
# Applies updated odds values to Sheet1 per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 5.5
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.95
$ws.Range("V2").Value = 1.23
$ws.Range("AA2").Value = 2.38
$ws.Range("AB2").Value = 1.53
$ws.Range("AG2").Value = 17
$ws.Range("AN2").Value = 11
$ws.Range("AO2").Value = 26
$ws.Range("G4").Value = 2.9
$ws.Range("I4").Value = 2.55
$ws.Range("L4").Value = 3.4
$ws.Range("AD4").Value = 13
$ws.Range("AK4").Value = 19
$ws.Range("AN4").Value = 6.5
$ws.Range("AO4").Value = 11
$ws.Range("G5").Value = 2.85
$ws.Range("I5").Value = 3.05
$ws.Range("J5").Value = 3.7
$ws.Range("K5").Value = 1.72
$ws.Range("L5").Value = 3.85
$ws.Range("AC5").Value = 5.9
$ws.Range("AE5").Value = 11.25
$ws.Range("AG5").Value = 35
$ws.Range("AJ5").Value = 5.1
$ws.Range("AN5").Value = 6.4
$ws.Range("AP5").Value = 11.5
$ws.Range("AR5").Value = 37
$ws.Range("AS5").Value = 55
$ws.Range("N11").Value = 9
$ws.Range("O11").Value = 1.36
$ws.Range("P11").Value = 3
$ws.Range("H12").Value = 3.3
$ws.Range("J12").Value = 3
$ws.Range("K12").Value = 2.05
$ws.Range("L12").Value = 4
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 9
$ws.Range("O12").Value = 1.36
$ws.Range("P12").Value = 3
$ws.Range("S12").Value = 2.2
$ws.Range("T12").Value = 1.65
$ws.Range("W12").Value = 4
$ws.Range("X12").Value = 1.22
$ws.Range("Y12").Value = 1.5
$ws.Range("Z12").Value = 2.5
$ws.Range("AA12").Value = 1.95
$ws.Range("AB12").Value = 1.8
$ws.Range("AC12").Value = 7
$ws.Range("AE12").Value = 9.5
$ws.Range("AH12").Value = 34
$ws.Range("AI12").Value = 8.5
$ws.Range("AK12").Value = 17
$ws.Range("AM12").Value = 351
$ws.Range("AN12").Value = 8.5
$ws.Range("AR12").Value = 29
$ws.Range("AS12").Value = 41
$ws.Range("N13").Value = 13
$ws.Range("S13").Value = 1.73
$ws.Range("T13").Value = 2.08
$ws.Range("W13").Value = 2.75
$ws.Range("X13").Value = 1.4
$ws.Range("G14").Value = 2.15
$ws.Range("I14").Value = 3.1
$ws.Range("J14").Value = 2.75
$ws.Range("K14").Value = 2.1
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 10
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 3.5
$ws.Range("S14").Value = 2
$ws.Range("T14").Value = 1.85
$ws.Range("W14").Value = 3.4
$ws.Range("X14").Value = 1.3
$ws.Range("Y14").Value = 1.4
$ws.Range("Z14").Value = 2.75
$ws.Range("AA14").Value = 1.8
$ws.Range("AB14").Value = 1.91
$ws.Range("AC14").Value = 8
$ws.Range("AE14").Value = 9.5
$ws.Range("AI14").Value = 10
$ws.Range("AN14").Value = 10
$ws.Range("AO14").Value = 17
$ws.Range("AQ14").Value = 34
